# api_function_complete_tracker.xlsx
# - added integration tests for Brazilian atlas, fixed up tests for Australian atlas

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Australian atlas (row 2) row got reset: clear the "X" status + the
#     "Unsure what to do" comment in col F, and give B2 the same (blank,
#     green) formatting as its neighbouring cells (C2:E2) ---
$ws.Range("C2").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B2").Value = ""

# --- Clear the generic "Unsure what to do" comment from col F for every
#     country row (2-12); the status is no longer applicable ---
$ws.Range("F2").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("F12").Value = ""

# --- Brazil (row 4): integration tests written now, so the "unable to
#     verify user details" note in E4 no longer applies; row also no
#     longer needs the taller custom row height ---
$ws.Range("E4").Value = ""
$ws.Rows.Item(4).AutoFit() | Out-Null

# --- Column F was sized for the old, longer comment text; shrink it back
#     down now that the comments are gone ---
$ws.Columns.Item(6).ColumnWidth = 14.6

# --- Selection moved to B2 (the cell that was just cleared) ---
$ws.Range("B2").Select() | Out-Null
